$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append, matching the diff: rows 34, 35 (duplicate of row 33),
# and row 36 (a new DecisionTree result).
$rows = @(
    @("LogisticRegression", "{'max_iter': 250}", 0.5011377861893516, 0.8325917431192661, 0.5094879356504962, 0.8287155963302753),
    @("LogisticRegression", "{'max_iter': 250}", 0.5011377861893516, 0.8325917431192661, 0.5094879356504962, 0.8287155963302753),
    @("DecisionTree", "{'max_depth': 30}", 0.05961969743960245, 0.9872935779816514, 3.443732803030968, 0.8844495412844037)
)

$startRow = 34
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
